$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questionnaire Results")

# Seed the new participant names first so shared-string indices are
# allocated in the same order as the source edit (Aimee, Sandra, Ada,
# Milica, Yves, then the new free-text activity string).
$ws.Range("B18").Value = "Aimee"
$ws.Range("B19").Value = "Sandra"
$ws.Range("B20").Value = "Ada"
$ws.Range("B21").Value = "Milica"
$ws.Range("B22").Value = "Yves"

# New respondents appended to the results table (rows 18-22 / IDs 15-19)
$ws.Range("A18").Value = 15
$ws.Range("C18").Value = "B"
$ws.Range("D18").Value = "A"
$ws.Range("E18").Value = "B"
$ws.Range("F18").Value = "Tipp B"
$ws.Range("I18").Value = "Tipp C"

$ws.Range("A19").Value = 16
$ws.Range("C19").Value = "B"
$ws.Range("D19").Value = "B"
$ws.Range("E19").Value = "B"
$ws.Range("F19").Value = "Tipp B"
$ws.Range("G19").Value = "Shopping or going to a museum"
$ws.Range("I19").Value = "Tipp A"

$ws.Range("A20").Value = 17
$ws.Range("C20").Value = "B"
$ws.Range("D20").Value = "C"
$ws.Range("E20").Value = "B"
$ws.Range("F20").Value = "Tipp C"
$ws.Range("G20").Value = "Shopping"
$ws.Range("I20").Value = "Tipp B"

$ws.Range("A21").Value = 18
$ws.Range("C21").Value = "A"
$ws.Range("D21").Value = "A"
$ws.Range("E21").Value = "A"
$ws.Range("F21").Value = "Tipp B"
$ws.Range("I21").Value = "Tipp B"

$ws.Range("A22").Value = 19
$ws.Range("C22").Value = "B"
$ws.Range("D22").Value = "A"
$ws.Range("E22").Value = "B"
$ws.Range("F22").Value = "Tipp C"
$ws.Range("G22").Value = "Shopping"
$ws.Range("I22").Value = "Tipp C"

# Touch the shared "Match" formula on the freshly populated rows so the
# cached results pick up the new F/I inputs (these cells held a cached
# blank-string result that otherwise survives straight recalculation).
for ($r = 18; $r -le 22; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.Formula = $cell.Formula
}

# Move the active selection to match the saved cursor position
$ws.Range("G29").Select()
